# Refresh crypto price (D) and 1h volume change (E) columns with the
# latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.345.52'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '3.620.51'
$ws.Range("E3").Value = '  +2.30%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = "'602.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = "'196.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  -1.02%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = "'0.213"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.72%  '
$ws.Range("D10").Value = "'0.646"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").Value = "'53.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.08%  '
$ws.Range("D12").Value = "'0.0000306"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("D13").Value = "'9.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.47%  '
$ws.Range("D14").Value = '4.195.90'
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("D15").Value = "'600.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.33%  '
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").Value = '70.513.87'
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("D18").Value = '3.635.01'
$ws.Range("E18").Value = '  +2.87%  '
$ws.Range("D19").Value = "'19.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = '  +1.63%  '
$ws.Range("D21").Value = "'0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("D22").Value = "'18.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.75%  '
$ws.Range("E23").Value = '  -1.79%  '
$ws.Range("D24").Value = "'102.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = "'4.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").Value = "'3.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.58%  '
$ws.Range("D27").Value = "'10.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.96%  '
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("D29").Value = "'33.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").Value = "'4.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.16%  '
$ws.Range("D31").Value = "'7.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.00%  '
$ws.Range("D32").Value = "'12.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.75%  '
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("D34").Value = "'63.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").Value = '0.0₃0884'
$ws.Range("E35").Value = '  +2.69%  '
$ws.Range("D36").Value = '3.919.83'
$ws.Range("E36").Value = '  +5.09%  '
$ws.Range("D37").Value = "'536.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.56%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("E39").Value = '  +1.10%  '
$ws.Range("D40").Value = "'36.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("D42").Value = "'3.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.04%  '
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").Value = "'0.0460"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.65%  '
$ws.Range("D45").Value = "'3.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.67%  '
$ws.Range("E46").Value = '  +1.09%  '
$ws.Range("D47").Value = "'0.140"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("E48").Value = '  -0.59%  '
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").Value = "'0.000251"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("E51").Value = '  +1.28%  '
